# Add the "TermExtList" worksheet with header + data, then format header bold.
$wb = $excel.ActiveWorkbook

$newSheet = $wb.Worksheets.Add()
$newSheet.Name = "TermExtList"

$headers = @("Publ. No", "Term Ext. [days]", "Disclaimer/Date", "Appl ID")
for ($c = 1; $c -le 4; $c++) {
    $newSheet.Cells.Item(1, $c).Value = $headers[$c - 1]
}
$newSheet.Range("A1:D1").Font.Bold = $true

$data = @(
    @("US20190107969A1", "83", "Terminal Disclaimer Filed2020-09-14 00:00:00", "16204798"),
    @("10272200", "no extension", "no disclaimer", "15114834"),
    @("US20210236729A1", "no extension", "no disclaimer", "17161528"),
    @("US20210038163A1", "no extension", "no disclaimer", "16533470"),
    @("9974492", "no extension", "no disclaimer", "15255909"),
    @("10185513", "167", "no disclaimer", "15256137")
)

for ($r = 0; $r -lt $data.Count; $r++) {
    for ($c = 0; $c -lt 4; $c++) {
        $newSheet.Cells.Item($r + 2, $c + 1).Value = $data[$r][$c]
    }
}

$newSheet.Columns("A:D").ColumnWidth = 30

# Move the new sheet to the end (after Tabelle1)
$newSheet.Move($null, $wb.Worksheets.Item($wb.Worksheets.Count))
